$d = $word.ActiveDocument

# --- 1. Insert the new paragraph right after the "Nedan presenteras fynd..." paragraph ---
$anchorText = "Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC, Chain of Custody, Controlled Wood och PEFC."
$newParaText = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text.Contains($anchorText)) {
        $par.Range.InsertParagraphAfter()
        $d.Paragraphs.Item($i + 1).Range.Text = $newParaText
        break
    }
}

# --- 2. Remove the two trailing empty paragraphs and the duplicate "Vi förväntar..." paragraph
#        that used to sit at the end of the document, right before the sectPr. Find the LAST
#        occurrence of the "artskyddsförordningen" comment paragraph (the one near the very end)
#        and delete everything between it and the end of the document body. ---
$lastText = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"

$count = $d.Paragraphs.Count
$lastMatch = -1
for ($i = 1; $i -le $count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text.Contains($lastText)) {
        $lastMatch = $i
    }
}

if ($lastMatch -ge 1) {
    $count = $d.Paragraphs.Count
    $startPara = $lastMatch + 1
    $endPara = $count
    if ($startPara -le $endPara) {
        $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
        $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
        $rg = $d.Range($rangeStart, $rangeEnd)
        $rg.Delete()
    }
}

# --- 3. Update the date shown in the title-page header from 2023-11-13 to 2023-11-14 ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$hdrFull = $hdr.Range.Text
$dateIdx = $hdrFull.IndexOf("2023-11-13")
if ($dateIdx -ge 0) {
    $dateRange = $hdr.Range.Duplicate
    $dateRange.SetRange($dateIdx, $dateIdx + 10)
    $dateRange.Text = "2023-11-14"
}
